# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recomputed upstream and re-written
# for this save file. Apply the recomputed values to the matching rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 3
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2
